$d = $word.ActiveDocument

# Find the Date paragraph (ends in "...October 8, 2025") and collapse to its end.
$find = $d.Content
$find.Find.Execute("2025", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$r = $find.Duplicate
$r.Collapse(0)

$datePara = $r.Paragraphs(1)

# Insert a new paragraph right after the Date paragraph: the Generative AI disclaimer.
$datePara.Range.InsertParagraphAfter() | Out-Null
$p1 = $datePara.Next()
$p1.Style = "FirstParagraph"
$p1.Range.Text = "Generative AI Disclaimer: AI was used to assist with organzing and arranging data pipelines to perform models."

# Insert another paragraph after the disclaimer: the GitHub repository link.
$p1.Range.InsertParagraphAfter() | Out-Null
$p2 = $p1.Next()
$p2.Style = "BodyText"
$p2.Range.Text = "Github Repository Link: https://github.com/met-ad-688/assignment-04-ryanmt1998.git"
